$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E in this sheet are text-typed (inlineStr) even when the
# content looks numeric (e.g. "246.57") or like a percentage (e.g. "+0.82%").
# A plain Range.Value assignment lets Excel auto-coerce numeric-looking
# strings into the Number type (and introduces float rounding artifacts), so
# force the cell to Text format first, then reset the style back to Normal
# afterwards so no stray style index is left on the cell (matches original).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '35.324.55'
Set-TextValue "E2" '  +0.82%  '
Set-TextValue "D3" '1.884.67'
Set-TextValue "E3" '  -0.20%  '
Set-TextValue "E4" '  -0.74%  '
Set-TextValue "D5" '246.57'
Set-TextValue "E5" '  -1.68%  '
Set-TextValue "E6" '  -1.03%  '
Set-TextValue "E7" '  -0.81%  '
Set-TextValue "D8" '43.38'
Set-TextValue "E8" '  +5.72%  '
Set-TextValue "E9" '  +1.55%  '
Set-TextValue "D10" '53.30'
Set-TextValue "E10" '  +2.19%  '
Set-TextValue "D11" '0.0750'
Set-TextValue "E11" '  +1.00%  '
Set-TextValue "D12" '0.0979'
Set-TextValue "E12" '  +0.42%  '
Set-TextValue "D13" '13.48'
Set-TextValue "E13" '  +4.68%  '
Set-TextValue "D14" '2.158.07'
Set-TextValue "E14" '  -0.29%  '
Set-TextValue "D15" '0.771'
Set-TextValue "E15" '  +7.03%  '
Set-TextValue "E16" '  +0.31%  '
Set-TextValue "D17" '1.903.35'
Set-TextValue "E17" '  +0.83%  '
Set-TextValue "D18" '35.318.67'
Set-TextValue "E18" '  +0.91%  '
Set-TextValue "D19" '73.99'
Set-TextValue "E19" '  +0.87%  '
Set-TextValue "E20" '  +0.31%  '
Set-TextValue "D21" '244.74'
Set-TextValue "E21" '  -2.39%  '
Set-TextValue "D22" '12.82'
Set-TextValue "E22" '  -0.22%  '
Set-TextValue "D23" '5.19'
Set-TextValue "E23" '  +4.68%  '
Set-TextValue "E24" '  +8.99%  '
Set-TextValue "E25" '  -0.72%  '
Set-TextValue "D26" '2.17'
Set-TextValue "E26" '  -2.47%  '
Set-TextValue "D27" '165.23'
Set-TextValue "E27" '  -1.10%  '
Set-TextValue "D28" '8.64'
Set-TextValue "E28" '  +2.00%  '
Set-TextValue "D29" '18.29'
Set-TextValue "E29" '  +0.08%  '
Set-TextValue "E30" '  +0.23%  '
Set-TextValue "E31" '  +0.83%  '
Set-TextValue "D32" '0.0595'
Set-TextValue "E32" '  +1.52%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D33" '1.89'
Set-TextValue "E33" '  -2.26%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D34" '4.20'
Set-TextValue "E34" '  +0.37%  '
Set-TextValue "E35" '  -0.75%  '
Set-TextValue "D36" '1.45'
Set-TextValue "E36" '  -4.64%  '
Set-TextValue "D37" '0.853'
Set-TextValue "E37" '  +1.69%  '
Set-TextValue "E38" '  -1.82%  '
Set-TextValue "D39" '0.0733'
Set-TextValue "E39" '  +10.59%  '
Set-TextValue "D40" '17.53'
Set-TextValue "E40" '  +0.24%  '
Set-TextValue "E41" '  +3.28%  '
Set-TextValue "E42" '  -0.77%  '
Set-TextValue "E43" '  -0.83%  '
Set-TextValue "E44" '  +1.64%  '
Set-TextValue "D45" '1.311.67'
Set-TextValue "E45" '  +1.36%  '
Set-TextValue "D46" '0.0803'
Set-TextValue "E46" '  +5.11%  '
Set-TextValue "E47" '  -1.13%  '
Set-TextValue "E48" '  -0.20%  '
Set-TextValue "D49" '11.92'
Set-TextValue "E49" '  -1.52%  '
Set-TextValue "E50" '  -2.07%  '
Set-TextValue "D51" '42.53'
Set-TextValue "E51" '  +1.25%  '
